# Add a new column name "pl_insol" (solar flux) by inserting a new row
# above "pl_temp_k" (habit code) in the single-column list on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting pl_temp_k and everything below it down
$ws.Rows("15:15").Insert()

# Populate the newly inserted cell
$ws.Cells.Item(15, 1).Value = "pl_insol"

# Mirror the resulting view/selection state (scroll + select the next empty row)
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$ws.Rows("28:28").Select() | Out-Null
